$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping-estado")

$map = @{
    "malo"       = "Malo"
    "bueno"      = "Bueno"
    "deficiente" = "Deficiente"
    "ruinoso"    = "Ruinoso"
}

for ($r = 1; $r -le 4; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
